# Generate Report for Handback
# Re-running the localization-status report updated the status of the
# "ad764da2-3354-4da9-831f-61edae4468bf" entry (row 3 on every sheet) from
# "Ready for handoff" to "Handed back: in sync with en-US", refreshed the
# handback datetime, and cleared the stale error detail message.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusHandedBack
$wsOverview.Range("F3").Value = $statusHandedBack

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusHandedBack
$wsZhCn.Range("K3").Value = "2016-08-13 07:00:57"
# Clear the stale error message but keep the cell as an (empty) text cell,
# matching the other blank cells in the table (e.g. L3/N3) instead of
# deleting it outright.
$wsZhCn.Range("P3").Value = "'"
$wsZhCn.Range("P3").Style = "Normal"

# --- de-de sheet ------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusHandedBack
$wsDeDe.Range("K3").Value = "2016-08-13 07:01:11"
$wsDeDe.Range("P3").Value = "'"
$wsDeDe.Range("P3").Style = "Normal"

# Column P (Error Detail) now holds only short/empty values; shrink it back
# down from the wide (40) autofit width used for the long error message.
$wsZhCn.Columns.Item(16).ColumnWidth = 12.83
$wsDeDe.Columns.Item(16).ColumnWidth = 12.83
